# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price ("D") and Volume(1h) ("E") columns are plain text cells in this sheet,
# so numeric-looking price strings are written with a leading apostrophe to
# keep them stored as text (matching the original inlineStr text cells)
# instead of being auto-converted to numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.690.28"
$ws.Range("E2").Value = "  -2.15%  "

$ws.Range("D3").Value = "2.294.95"
$ws.Range("E3").Value = "  -3.44%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'301.41"
$ws.Range("E5").Value = "  -2.82%  "

$ws.Range("D6").Value = "'98.16"
$ws.Range("E6").Value = "  -6.70%  "

$ws.Range("D7").Value = "'0.506"
$ws.Range("E7").Value = "  -1.17%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "'0.504"
$ws.Range("E9").Value = "  -3.32%  "

$ws.Range("D10").Value = "'34.62"
$ws.Range("E10").Value = "  -4.46%  "

$ws.Range("D11").Value = "'51.23"
$ws.Range("E11").Value = "  -4.16%  "

$ws.Range("D12").Value = "'0.0792"
$ws.Range("E12").Value = "  -2.78%  "

$ws.Range("E13").Value = "  +0.41%  "

$ws.Range("D14").Value = "'6.73"
$ws.Range("E14").Value = "  -4.26%  "

$ws.Range("E15").Value = "  -3.20%  "

$ws.Range("D16").Value = "'15.43"
$ws.Range("E16").Value = "  -1.65%  "

$ws.Range("D17").Value = "2.289.32"
$ws.Range("E17").Value = "  -3.81%  "

$ws.Range("D18").Value = "'0.793"
$ws.Range("E18").Value = "  -2.65%  "

$ws.Range("D19").Value = "42.709.58"

$ws.Range("D20").Value = "'11.63"
$ws.Range("E20").Value = "  -2.81%  "

$ws.Range("D21").Value = "0.0₃0897"
$ws.Range("E21").Value = "  -2.47%  "

$ws.Range("D22").Value = "'6.03"
$ws.Range("E22").Value = "  -4.62%  "

$ws.Range("D23").Value = "'67.22"
$ws.Range("E23").Value = "  -1.76%  "

$ws.Range("D24").Value = "'235.17"
$ws.Range("E24").Value = "  -2.58%  "

$ws.Range("D25").Value = "'1.95"
$ws.Range("E25").Value = "  -5.08%  "

$ws.Range("E26").Value = "  -4.44%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").Value = "'24.64"
$ws.Range("E28").Value = "  -4.48%  "

$ws.Range("E29").Value = "  -4.78%  "

$ws.Range("D30").Value = "'34.22"
$ws.Range("E30").Value = "  -6.92%  "

$ws.Range("D31").Value = "'164.66"
$ws.Range("E31").Value = "  +2.06%  "

$ws.Range("D32").Value = "'9.11"
$ws.Range("E32").Value = "  -4.69%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").Value = "'4.99"
$ws.Range("E34").Value = "  -5.47%  "

$ws.Range("E35").Value = "  -4.92%  "

$ws.Range("D36").Value = "'0.0701"
$ws.Range("E36").Value = "  -5.68%  "

$ws.Range("D37").Value = "'4.37"
$ws.Range("E37").Value = "  -6.32%  "

$ws.Range("D38").Value = "'2.84"
$ws.Range("E38").Value = "  -8.49%  "

$ws.Range("D39").Value = "'16.24"
$ws.Range("E39").Value = "  -11.45%  "

$ws.Range("D40").Value = "'1.79"
$ws.Range("E40").Value = "  -7.85%  "

$ws.Range("E41").Value = "  -4.94%  "

$ws.Range("D42").Value = "'0.110"
$ws.Range("E42").Value = "  -3.00%  "

$ws.Range("D43").Value = "'2.42"
$ws.Range("E43").Value = "  -7.18%  "

$ws.Range("D44").Value = "1.973.13"
$ws.Range("E44").Value = "  -3.07%  "

$ws.Range("E45").Value = "  -2.21%  "

$ws.Range("D46").Value = "'18.22"
$ws.Range("E46").Value = "  -8.00%  "

$ws.Range("D47").Value = "'9.81"
$ws.Range("E47").Value = "  -7.14%  "

$ws.Range("D48").Value = "'2.87"
$ws.Range("E48").Value = "  -9.44%  "

$ws.Range("D49").Value = "'4.75"
$ws.Range("E49").Value = "  +0.00%  "

# Row 50 / 51: coin order swapped (MultiversX now ranks above RocketPoolETH)
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "'53.50"
$ws.Range("E50").Value = "  -7.69%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.528.43"
$ws.Range("E51").Value = "  -2.94%  "
